# Added 4wk low sales check - updates forecast numbers (MyForecast, Seasonality
# Index) on the "Forecast Comparison" sheet and the derived roll-up totals on
# the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D) and Seasonality Index (L) ---------

$forecastUpdates = @(
    @{ Row = 2;  D = 11; L = 1.13 },
    @{ Row = 3;  D = 15; L = 0.89 },
    @{ Row = 4;  D = 16; L = 1 },
    @{ Row = 5;  D = 15; L = 0.9 },
    @{ Row = 6;  D = 13; L = 0.99 },
    @{ Row = 7;  D = 12; L = 0.93 },
    @{ Row = 8;  D = 13; L = 0.8 },
    @{ Row = 9;  D = 15; L = 1.11 },
    @{ Row = 10; D = 16; L = 0.99 },
    @{ Row = 11; D = 14; L = 1.11 },
    @{ Row = 12; D = 12; L = 0.98 },
    @{ Row = 13; D = 12 },
    @{ Row = 14; D = 14 },
    @{ Row = 15; D = 15; L = 0.87 },
    @{ Row = 16; D = 15; L = 1.08 },
    @{ Row = 17; D = 14; L = 0.97 }
)

foreach ($u in $forecastUpdates) {
    $wsForecast.Cells.Item($u.Row, 4).Value = $u.D
    if ($u.ContainsKey("L")) {
        $wsForecast.Cells.Item($u.Row, 12).Value = $u.L
    }
}

# --- Summary: derived totals stored as text in column B --------------------
# Leading apostrophe keeps these numeric-looking values stored as text,
# matching the existing (inline-string) cell type on this sheet.

$wsSummary.Range("B9").Value  = "'222"
$wsSummary.Range("B10").Value = "'110"
$wsSummary.Range("B11").Value = "'57"
$wsSummary.Range("B12").Value = "'16"
$wsSummary.Range("B14").Value = "'11"
